# Insert a new record (account 004584517 / CAIO / 10000) as the new
# third row of the "Export" sheet, right after the THIAGO row and
# immediately before the ALPHASITIO (005305448) row, shifting every
# row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Push row 3 (and everything after) down to make room for the new row.
$ws.Rows.Item(3).Insert()

# Column A holds account numbers with significant leading zeros, so it
# must stay text (otherwise "004584517" would be read as the number
# 4584517). Force the text format before assigning the value.
$acct = $ws.Cells.Item(3, 1)
$acct.NumberFormat = "@"
$acct.Value = "004584517"

$ws.Cells.Item(3, 2).Value = "CAIO"
$ws.Cells.Item(3, 3).Value = 10000
